# Add some other MCPD fields (DONORCODE, DONORNUMB, DUPLSITE, STORAGE) to the
# accessions worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ESP058:CIAM81001 accession)
$ws.Range("AC2").Value = "ESP058"
$ws.Range("AD2").Value = "ecu1"
$ws.Range("AF2").Value = "algunlao; otro lado"
$ws.Range("AG2").Value = 13

# Row 3 (ESP004:BGE005836 accession)
$ws.Range("AD3").Value = "ecu2"
$ws.Range("AG3").Value = 12

# Row 4 (ESP058:CIAM81002 accession)
$ws.Range("AD4").Value = "ecu3"
$ws.Range("AG4").Value = 12

# Row 5 (ESP004:BGE005837 accession)
$ws.Range("AD5").Value = "ecu4"
$ws.Range("AG5").Value = 13

# Row 6 (ESP058:CIAM81003 accession)
$ws.Range("AD6").Value = "ecu5"
$ws.Range("AG6").Value = 12

# Row 7 (ESP004:BGE005838 accession)
$ws.Range("AD7").Value = "ecu6"

# Reflect the author's final selection/view state (active cell AG2).
[void]$ws.Range("AG2").Select()
